$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) ODI Batting: drop the (empty) INNING_NUMBER cells for rows that
#    have no value - matches rows 2, 5 and 6 (B2, B5, B6) in the diff.
# ------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B5").ClearContents()
$batting.Range("B6").ClearContents()

# ------------------------------------------------------------------
# 2) Add the new "ODI Batting Extra" sheet after "ODI Bowling".
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row (bold, centered, bordered - same look as the other sheets).
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$headerRange = $extra.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data rows (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH).
$data = @(
    @("4472", 9,    "",  "",  "",       "NO"),
    @("4473", 9,    "2", "0", "12.55%", "NO"),
    @("4476", 9,    "1", "0", "3.61%",  "NO"),
    @("4599", "",   "",  "",  "",       "NO"),
    @("4602", 9,    "",  "",  "",       "NO"),
    @("4609", 10,   "2", "0", "13.64%", "NO"),
    @("4613", 10,   "0", "0", "0.81%",  "NO"),
    @("4618", 10,   "0", "0", "1.16%",  "NO"),
    @("4619", "",   "",  "",  "",       "NO")
)

$row = 2
foreach ($rec in $data) {
    # A: MATCH_CODE - numeric-looking but kept as text.
    $codeCell = $extra.Cells.Item($row, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $rec[0]

    # B: BATTING_POSITION - a real number when present, an (empty) cell otherwise.
    $posCell = $extra.Cells.Item($row, 2)
    if ($rec[1] -ne "") {
        $posCell.Value = $rec[1]
    } else {
        $posCell.Borders.LineStyle = -4142
    }

    # C: NUM_4, D: NUM_6 - numeric-looking but kept as text; (empty) cell if blank.
    foreach ($colIdx in 3, 4) {
        $cell = $extra.Cells.Item($row, $colIdx)
        $val = $rec[$colIdx - 1]
        if ($val -ne "") {
            $cell.NumberFormat = "@"
            $cell.Value = $val
        } else {
            $cell.Borders.LineStyle = -4142
        }
    }

    # E: PERCENT_RUNS_OF_TOTAL - textual percentage, not a numeric percent; (empty) cell if blank.
    $pctCell = $extra.Cells.Item($row, 5)
    $pctVal = $rec[4]
    if ($pctVal -ne "") {
        $pctCell.NumberFormat = "@"
        $pctCell.Value = $pctVal
    } else {
        $pctCell.Borders.LineStyle = -4142
    }

    # F: MAN_OF_MATCH - plain text.
    $extra.Cells.Item($row, 6).Value = $rec[5]

    $row++
}
